# Fruta / hortaliza, semanal
# Insert a new weekly record for "Chirimoya" (Vega Modelo de Temuco) at
# row 100, pushing the existing rows 100-140 down to 101-141.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 100 - Excel shifts rows
# 100..140 down to 101..141 and carries formatting (incl. the date
# style on column D) down from row 99.
$ws.Rows.Item(100).Insert()

# Populate the newly inserted row 100 with this week's Chirimoya data.
$ws.Cells.Item(100, 1).Value = 10
$ws.Cells.Item(100, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(100, 3).Value = "La Araucanía"
$ws.Cells.Item(100, 4).Value = 44809
$ws.Cells.Item(100, 5).Value = 9
$ws.Cells.Item(100, 6).Value = "Fruta"
$ws.Cells.Item(100, 7).Value = 100107
$ws.Cells.Item(100, 8).Value = "Otros"
$ws.Cells.Item(100, 9).Value = 100107002
$ws.Cells.Item(100, 10).Value = "Chirimoya"
$ws.Cells.Item(100, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(100, 12).Value = "Primera"
$ws.Cells.Item(100, 13).Value = 95
$ws.Cells.Item(100, 14).Value = 3500
$ws.Cells.Item(100, 15).Value = 3500
$ws.Cells.Item(100, 16).Value = 3500
$ws.Cells.Item(100, 17).Value = "$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(100, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(100, 19).Value = 3500
$ws.Cells.Item(100, 20).Value = 1
